# Generate Report for Handback
#
# The handback transform failed for the "f08df354-8bb7-40cd-9e99-592ddd447fe9"
# file in both locales: the status flips from "Ready for handoff" to
# "Handback transform failed", and an explanatory error message is written
# into the "Error Detail" column (K) of each locale's detail sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the f08df354-... file; B=zh-cn status, C=de-de status.
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# zh-cn detail sheet: row 3 is the f08df354-... file; C = Status column.
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "Handback file name: l3dtazol.jbg is different with handoff file name: f08df354-8bb7-40cd-9e99-592ddd447fe9.17bbac4f90ee89186996ec41f7712f755d72aaf4.zh-cn."

# de-de detail sheet: row 3 is the f08df354-... file; C = Status column.
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "Handback file name: l3dtazol.jbg is different with handoff file name: f08df354-8bb7-40cd-9e99-592ddd447fe9.17bbac4f90ee89186996ec41f7712f755d72aaf4.de-de."
